# The workbook used to ship two sheets: a computed "Total Report" summary
# sheet and the raw "Form Responses 1" data sheet. This edit drops the
# (now unused/unneeded) "Total Report" sheet, leaving only the form
# responses. Excel automatically re-derives the dependent bits on save:
#   - the shared-string table sheds the two strings ("ยอดขาย(บาท)",
#     "รวม") that were only referenced from the deleted sheet;
#   - the remaining (now only/active) worksheet's view picks up
#     tabSelected="1".
$wb = $excel.ActiveWorkbook
[void]$wb.Worksheets.Item("Total Report").Delete()

# The default "Normal" cell style had been persisted under its Thai
# localized display name ("ปกติ"). Deleting it collapses it back to the
# engine's canonical built-in "Normal" style.
[void]$wb.Styles.Item(1).Delete()
